# Site updated: 2017-08-15 15:02:19
#
# Header row (row 1) gets a taller, vertically-centered look, row 2 grows a
# touch too, and the saved view no longer remembers a mid-sheet scroll
# position / stray active-cell offset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1: taller header row, vertically centered content ---------------
$ws.Rows.Item(1).RowHeight = 23.25
$ws.Range("B1:D1").VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter

# --- Row 2: slightly taller ------------------------------------------------
$ws.Rows.Item(2).RowHeight = 16.5

# --- Reset the saved view: scroll back to the top, clear the stray active
#     cell offset so the selection is just the used range ------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A1:D17").Select()
$excel.ActiveCell.Worksheet.Range("A1").Select()
$ws.Range("A1:D17").Select()
